# Apply the changes described in the diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change: O1
$ws.Range("O1").Value = "F1 train"

# Column O (Validation -> F1 train) numeric updates, rows 2-15
$ws.Range("O2").Value = 0.927536231884058
$ws.Range("O3").Value = 0.9855072463768116
$ws.Range("O4").Value = 0.9855072463768116
$ws.Range("O5").Value = 0.5660377358490566
$ws.Range("O6").Value = 0.6376811594202898
$ws.Range("O7").Value = 0.9855072463768116
$ws.Range("O8").Value = 1
$ws.Range("O9").Value = 0.9565217391304348
$ws.Range("O10").Value = 0.8125
$ws.Range("O11").Value = 0.5333333333333333
$ws.Range("O12").Value = 0.9411764705882353
$ws.Range("O13").Value = 0.9855072463768116
$ws.Range("O14").Value = 1
$ws.Range("O15").Value = 1

# Row 16 updates (MLP / Free technique row)
$ws.Range("C16").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 0.5
$ws.Range("J16").Value = 0.5833333333333334
$ws.Range("K16").Value = 0.7777777777777778
$ws.Range("L16").Value = 0.4666666666666667
$ws.Range("M16").Value = 0.2727272727272727
$ws.Range("N16").Value = 0.7777777777777778
$ws.Range("O16").Value = 0.7848101265822784
